$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '39.846.75'
$ws.Cells.Item(2, 5).Value = '  +0.45%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.223.66'
$ws.Cells.Item(3, 5).Value = '  +0.64%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.02%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''291.65'
$ws.Cells.Item(5, 5).Value = '  -0.05%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''87.04'
$ws.Cells.Item(6, 5).Value = '  +0.88%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''0.513'
$ws.Cells.Item(7, 5).Value = '  -0.24%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.03%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -0.62%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''30.46'
$ws.Cells.Item(10, 5).Value = '  +0.48%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''0.0781'
$ws.Cells.Item(11, 5).Value = '  -0.48%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''49.93'
$ws.Cells.Item(12, 5).Value = '  +5.53%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +2.72%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''6.44'
$ws.Cells.Item(14, 5).Value = '  +1.79%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '2.570.14'
$ws.Cells.Item(15, 5).Value = '  +0.85%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '''13.82'
$ws.Cells.Item(16, 5).Value = '  -1.49%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '2.228.90'
$ws.Cells.Item(17, 5).Value = '  +0.81%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '''0.732'
$ws.Cells.Item(18, 5).Value = '  +0.74%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '39.796.70'
$ws.Cells.Item(19, 5).Value = '  +0.42%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '0.0₃0886'
$ws.Cells.Item(20, 5).Value = '  +0.80%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''11.07'
$ws.Cells.Item(21, 5).Value = '  -2.75%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''5.75'
$ws.Cells.Item(22, 5).Value = '  -0.67%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''65.78'
$ws.Cells.Item(23, 5).Value = '  +0.18%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''237.11'
$ws.Cells.Item(24, 5).Value = '  +0.66%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -0.09%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''2.45'
$ws.Cells.Item(26, 5).Value = '  -0.19%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +0.06%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''23.09'
$ws.Cells.Item(28, 5).Value = '  +1.76%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '''9.24'
$ws.Cells.Item(29, 5).Value = '  -0.18%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -6.88%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''156.78'
$ws.Cells.Item(31, 5).Value = '  +3.17%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''31.90'
$ws.Cells.Item(32, 5).Value = '  -2.37%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -0.03%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +1.12%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +7.26%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -0.33%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -1.66%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -0.24%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''1.75'
$ws.Cells.Item(39, 5).Value = '  +3.66%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''0.0990'
$ws.Cells.Item(40, 5).Value = '  +0.44%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''15.28'
$ws.Cells.Item(41, 5).Value = '  -4.08%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '2.108.43'
$ws.Cells.Item(42, 5).Value = '  +1.66%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''3.72'
$ws.Cells.Item(43, 5).Value = '  -1.70%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'EnergySwap'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(44, 4).Value = '''18.16'
$ws.Cells.Item(44, 5).Value = '  +3.12%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'VeChain'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(45, 4).Value = '''0.0271'
$ws.Cells.Item(45, 5).Value = '  +1.45%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''9.97'
$ws.Cells.Item(46, 5).Value = '  +0.09%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -7.19%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''2.72'
$ws.Cells.Item(48, 5).Value = '  +4.46%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '2.442.06'
$ws.Cells.Item(49, 5).Value = '  +0.92%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '''1.46'
$ws.Cells.Item(50, 5).Value = '  +2.38%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +2.71%  '
